# Update auto calibration 'all' to use meta y factor
# Rename the "Quantity Type" column header to "Units" on the Parameters
# and State Variables sheets, and update the active sheet/selection to
# match what was left selected in Excel (Parameters!N12).

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("Parameters")
$wsParams.Range("B1").Value = "Units"
$wsParams.Range("B5").Value = "Units"
$wsParams.Range("B9").Value = "Units"
$wsParams.Range("B13").Value = "Units"
$wsParams.Range("B17").Value = "Units"

$wsState = $wb.Worksheets.Item("State Variables")
$wsState.Range("B1").Value = "Units"
$wsState.Range("B5").Value = "Units"
$wsState.Range("B9").Value = "Units"

$wsParams.Activate()
$wsParams.Range("N12").Select()
